$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("P4").Value = 49
$ws.Range("Q4").Value = 4.9
$ws.Range("R4").Value = 24.5
$ws.Range("S4").Value = 59
$ws.Range("T4").Value = 2940

# Row 5
$ws.Range("P5").Value = 49.5
$ws.Range("Q5").Value = 4.8
$ws.Range("R5").Value = 24
$ws.Range("S5").Value = 59
$ws.Range("T5").Value = 2880

# Row 6
$ws.Range("P6").Value = 47.5
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 9.5
$ws.Range("S6").Value = 57
$ws.Range("T6").Value = 1140

# Row 7
$ws.Range("P7").Value = 47.5
$ws.Range("Q7").Value = 1.9
$ws.Range("R7").Value = 9.5
$ws.Range("S7").Value = 57
$ws.Range("T7").Value = 1140

# Row 8
$ws.Range("P8").Value = 47.1311
$ws.Range("Q8").Value = 1.4
$ws.Range("R8").Value = 7
$ws.Range("S8").Value = 57
$ws.Range("T8").Value = 840

# Row 9
$ws.Range("P9").Value = 45.9016
$ws.Range("Q9").Value = 0.8
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 55
$ws.Range("T9").Value = 480

# Row 10
$ws.Range("P10").Value = 33.1967
$ws.Range("Q10").Value = 3.15
$ws.Range("R10").Value = 15.75
$ws.Range("S10").Value = 40
$ws.Range("T10").Value = 1890

# Row 11 (Q11, R11, T11 unchanged)
$ws.Range("P11").Value = 44.0574
$ws.Range("S11").Value = 53

# Row 12
$ws.Range("P12").Value = 47.541
$ws.Range("Q12").Value = 1.6
$ws.Range("R12").Value = 8
$ws.Range("S12").Value = 57
$ws.Range("T12").Value = 960

# Row 13
$ws.Range("P13").Value = 46.3115
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 9.5
$ws.Range("S13").Value = 56
$ws.Range("T13").Value = 1140

# Row 14
$ws.Range("P14").Value = 47.541
$ws.Range("Q14").Value = 1.6
$ws.Range("R14").Value = 8
$ws.Range("S14").Value = 57
$ws.Range("T14").Value = 960

# Row 15
$ws.Range("P15").Value = 37.3864
$ws.Range("Q15").Value = 5.45
$ws.Range("R15").Value = 27.25
$ws.Range("S15").Value = 45
$ws.Range("T15").Value = 3270
